# Scheduled market-data refresh for the leve-profit workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# on the rows whose backing market snapshot changed.
$wb = $excel.ActiveWorkbook

# ALC!62 - The Mustache Suits Him / Enchanted Mythrite Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 9043.143
$ws.Range("I62").Value = 9500
$ws.Range("J62").Value = 8860.4
$ws.Range("K62").Value = 9500
$ws.Range("L62").Value = 8860.4
$ws.Range("M62").Value = -8876
$ws.Range("N62").Value = -10108.4

# ALC!65 - Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 9043.143
$ws.Range("I65").Value = 9500
$ws.Range("J65").Value = 8860.4
$ws.Range("K65").Value = 47500
$ws.Range("L65").Value = 44302
$ws.Range("M65").Value = -44380
$ws.Range("N65").Value = -50542

# ALC!112 - Making Ends Meet / Superior Spiritbond Potion
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4311393
$ws.Range("J112").Value = 1088.8518
$ws.Range("L112").Value = 3266.5554
$ws.Range("N112").Value = -5482.555399999999

# ALC!116 - Growing Up / Growth Formula Kappa
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 12833293
$ws.Range("I116").Value = 25661150
$ws.Range("J116").Value = 5436.5454
$ws.Range("K116").Value = 25661150
$ws.Range("L116").Value = 5436.5454
$ws.Range("M116").Value = -25657708
$ws.Range("N116").Value = -12320.5454

# ALC!129 - Practical Command / Commanding Craftsman's Draught
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 250917.75
$ws.Range("J129").Value = 271240.9
$ws.Range("L129").Value = 813722.7000000001
$ws.Range("N129").Value = -823722.7000000001

# ARM!2 - Ain't Got No Ingots / Bronze Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1234.25
$ws.Range("I2").Value = 1224.2142
$ws.Range("K2").Value = 1224.2142
$ws.Range("M2").Value = -1111.2142

# ARM!32 - Ingot We Trust / Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13889.704
$ws.Range("I32").Value = 10277.156
$ws.Range("J32").Value = 23523.166
$ws.Range("K32").Value = 10277.156
$ws.Range("L32").Value = 23523.166
$ws.Range("M32").Value = -9990.156000000001
$ws.Range("N32").Value = -24097.166

# ARM!45 - Hollow Hallmarks / Mythril Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3119.2144
$ws.Range("I45").Value = 2747.8948
$ws.Range("J45").Value = 3903.111
$ws.Range("K45").Value = 2747.8948
$ws.Range("L45").Value = 3903.111
$ws.Range("M45").Value = -2370.8948
$ws.Range("N45").Value = -4657.111

# ARM!61 - Dealing with the Tough Stuff / Cobalt Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 17549184
$ws.Range("J61").Value = 5349.25
$ws.Range("L61").Value = 5349.25
$ws.Range("N61").Value = -5773.25

# ARM!116 - No Scope / Titanbronze Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1234.25
$ws.Range("I116").Value = 1224.2142
$ws.Range("K116").Value = 1224.2142
$ws.Range("M116").Value = 1069.7858

# ARM!132 - Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 11642074
$ws.Range("I132").Value = 13891226
$ws.Range("J132").Value = 75005.71000000001
$ws.Range("K132").Value = 41673678
$ws.Range("L132").Value = 225017.13
$ws.Range("M132").Value = -41671148
$ws.Range("N132").Value = -230077.13

# ARM!136 - Metal with Mettle / Cobalt Tungsten Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 17549184
$ws.Range("J136").Value = 5349.25
$ws.Range("L136").Value = 16047.75
$ws.Range("N136").Value = -21147.75

# BSM!3 - Hells Bells / Bronze Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1234.25
$ws.Range("I3").Value = 1224.2142
$ws.Range("K3").Value = 1224.2142
$ws.Range("M3").Value = -1110.2142

# BSM!94 - High Steal / High Steel Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1398.5
$ws.Range("I94").Value = 900.61536
$ws.Range("J94").Value = 2117.6667
$ws.Range("K94").Value = 900.61536
$ws.Range("L94").Value = 2117.6667
$ws.Range("M94").Value = -449.61536
$ws.Range("N94").Value = -3019.6667

# BSM!132 - Always Be Prepaired / Mountain Chromite Twinfangs
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 59897.25
$ws.Range("J132").Value = 59897.25
$ws.Range("L132").Value = 59897.25
$ws.Range("N132").Value = -70017.25

# BSM!133 - Paring Is Caring / Mountain Chromite Hatchet
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 45000
$ws.Range("J133").Value = 45000
$ws.Range("L133").Value = 45000
$ws.Range("N133").Value = -55120

# BSM!134 - Ruthenium Supremium / Ruthenium Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8380.200000000001
$ws.Range("I134").Value = 8831.923000000001
$ws.Range("J134").Value = 5444
$ws.Range("K134").Value = 26495.769
$ws.Range("L134").Value = 16332
$ws.Range("M134").Value = -23960.769
$ws.Range("N134").Value = -21402

# CRP!16 - Raise the Roof / Ash Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1150
$ws.Range("I16").Value = 1214
$ws.Range("K16").Value = 1214
$ws.Range("M16").Value = -927

# CRP!31 - Wall Not Found / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4488.4644
$ws.Range("I31").Value = 2226.1904
$ws.Range("J31").Value = 5845.8286
$ws.Range("K31").Value = 2226.1904
$ws.Range("L31").Value = 5845.8286
$ws.Range("M31").Value = -1931.1904
$ws.Range("N31").Value = -6435.8286

# CRP!34 - Armoires of the Rich and Famous / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4488.4644
$ws.Range("I34").Value = 2226.1904
$ws.Range("J34").Value = 5845.8286
$ws.Range("K34").Value = 2226.1904
$ws.Range("L34").Value = 5845.8286
$ws.Range("M34").Value = -2024.1904
$ws.Range("N34").Value = -6249.8286

# CRP!107 - Built to Last / White Oak Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1239.1613
$ws.Range("I107").Value = 776.7059
$ws.Range("J107").Value = 1800.7142
$ws.Range("K107").Value = 776.7059
$ws.Range("L107").Value = 1800.7142
$ws.Range("M107").Value = 1143.2941
$ws.Range("N107").Value = -5640.7142

# CRP!113 - Patient Patients / White Ash Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1150
$ws.Range("I113").Value = 1214
$ws.Range("K113").Value = 1214
$ws.Range("M113").Value = 956

# CRP!134 - Wood You Be Quiet / Ceiba Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 142858320
$ws.Range("I134").Value = 200001140
$ws.Range("K134").Value = 600003420
$ws.Range("M134").Value = -600000885

# CUL!64 - The Aroma of Faith / Baked Onion Soup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4260.6665
$ws.Range("I64").Value = 847.5
$ws.Range("K64").Value = 2542.5
$ws.Range("M64").Value = -2272.5

# CUL!67 - Soup's On (L) / Baked Onion Soup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 4260.6665
$ws.Range("I67").Value = 847.5
$ws.Range("K67").Value = 2542.5
$ws.Range("M67").Value = -1606.5

# CUL!86 - Let's Not Get Sappy / Birch Syrup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 340.36365
$ws.Range("I86").Value = 272.8
$ws.Range("J86").Value = 396.66666
$ws.Range("K86").Value = 818.4000000000001
$ws.Range("L86").Value = 1189.99998
$ws.Range("M86").Value = 367.5999999999999
$ws.Range("N86").Value = -3561.99998

# CUL!89 - Luxury Spillover (L) / Birch Syrup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 340.36365
$ws.Range("I89").Value = 272.8
$ws.Range("J89").Value = 396.66666
$ws.Range("K89").Value = 2455.2
$ws.Range("L89").Value = 3569.99994
$ws.Range("M89").Value = 3472.8
$ws.Range("N89").Value = -15425.99994

# CUL!131 - The Mountain Steeped / Tsai tou Vounou
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 723.54
$ws.Range("I131").Value = 395.75
$ws.Range("J131").Value = 768.23865
$ws.Range("K131").Value = 1187.25
$ws.Range("L131").Value = 2304.71595
$ws.Range("M131").Value = 3852.75
$ws.Range("N131").Value = -12384.71595

# GSM!4 - Arms for the Poor / Bone Brand
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

# WVR!2 - The Unmentionables / Hempen Underpants
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 111157520
$ws.Range("I2").Value = 250003170
$ws.Range("J2").Value = 81002.39999999999
$ws.Range("K2").Value = 250003170
$ws.Range("L2").Value = 81002.39999999999
$ws.Range("M2").Value = -250003058
$ws.Range("N2").Value = -81226.39999999999
